$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04941833333333333
$ws.Range("H2").Value = 0.148255
$ws.Range("I2").Value = 0.005167549122999764
$ws.Range("J2").Value = 0.005167549122999764
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05968133333333333
$ws.Range("N2").Value = 0.179044
$ws.Range("O2").Value = 0.02602747651633847
$ws.Range("P2").Value = 0.02602747651633848
$ws.Range("Q2").Value = 0.002949352024444444
$ws.Range("R2").Value = 0.02654416822
$ws.Range("S2").Value = 0.0001344982634459018
$ws.Range("T2").Value = 0.0001344982634459018
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04941833333333333
$ws.Range("H3").Value = 0.148255
$ws.Range("I3").Value = 0.005167549122999764
$ws.Range("J3").Value = 0.005167549122999764
$ws.Range("O3").Value = 0.144012433133819
$ws.Range("P3").Value = 0.144012433133819
$ws.Range("Q3").Value = 0.01631903734277778
$ws.Range("R3").Value = 0.146871336085
$ws.Range("S3").Value = 0.0007441913225417285
$ws.Range("T3").Value = 0.0007441913225417285
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.04941833333333333
$ws.Range("H4").Value = 0.148255
$ws.Range("I4").Value = 0.005167549122999764
$ws.Range("J4").Value = 0.005167549122999764
$ws.Range("O4").Value = 0.8299600903498424
$ws.Range("P4").Value = 0.8299600903498425
$ws.Range("Q4").Value = 0.09404847493166667
$ws.Range("R4").Value = 0.846436274385
$ws.Range("S4").Value = 0.004288859537012133
$ws.Range("T4").Value = 0.004288859537012133
$ws.Range("I5").Value = 0.806706161560336
$ws.Range("J5").Value = 0.806706161560336
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05968133333333333
$ws.Range("N5").Value = 0.179044
$ws.Range("O5").Value = 0.02602747651633847
$ws.Range("P5").Value = 0.02602747651633848
$ws.Range("Q5").Value = 0.4604233833288889
$ws.Range("R5").Value = 4.14381044996
$ws.Range("S5").Value = 0.02099652567559719
$ws.Range("T5").Value = 0.0209965256755972
$ws.Range("I6").Value = 0.806706161560336
$ws.Range("J6").Value = 0.806706161560336
$ws.Range("O6").Value = 0.144012433133819
$ws.Range("P6").Value = 0.144012433133819
$ws.Range("S6").Value = 0.1161757171503477
$ws.Range("T6").Value = 0.1161757171503477
$ws.Range("I7").Value = 0.806706161560336
$ws.Range("J7").Value = 0.806706161560336
$ws.Range("O7").Value = 0.8299600903498424
$ws.Range("P7").Value = 0.8299600903498425
$ws.Range("S7").Value = 0.6695339187343911
$ws.Range("T7").Value = 0.6695339187343912
$ws.Range("I8").Value = 0.1881262893166642
$ws.Range("J8").Value = 0.1881262893166643
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.05968133333333333
$ws.Range("N8").Value = 0.179044
$ws.Range("O8").Value = 0.02602747651633847
$ws.Range("P8").Value = 0.02602747651633848
$ws.Range("Q8").Value = 0.1073721098804444
$ws.Range("R8").Value = 0.9663489889239999
$ws.Range("S8").Value = 0.004896452577295375
$ws.Range("T8").Value = 0.004896452577295377
$ws.Range("I9").Value = 0.1881262893166642
$ws.Range("J9").Value = 0.1881262893166643
$ws.Range("O9").Value = 0.144012433133819
$ws.Range("P9").Value = 0.144012433133819
$ws.Range("S9").Value = 0.0270925246609296
$ws.Range("T9").Value = 0.0270925246609296
$ws.Range("I10").Value = 0.1881262893166642
$ws.Range("J10").Value = 0.1881262893166643
$ws.Range("O10").Value = 0.8299600903498424
$ws.Range("P10").Value = 0.8299600903498425
$ws.Range("S10").Value = 0.1561373120784392
$ws.Range("T10").Value = 0.1561373120784393
